$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Data PGK B.Makanan (β)")
$ws3 = $wb.Worksheets.Item("Data PGK B.Makanan (p)")

# --- Sheet "Data PGK B.Makanan (β)": add E (price relative) and F (computed) columns ---
$ws2.Range("E3").Value = 0.84126142999999998
$ws2.Range("E4").Value = 0.48830415999999999
$ws2.Range("E5").Value = 0.52143086999999999
$ws2.Range("E6").Value = 0.78539656000000002
$ws2.Range("E7").Value = 0.81992346000000005
$ws2.Range("E8").Value = 0.80531744999999999
$ws2.Range("E9").Value = 0.77272596000000005
$ws2.Range("E10").Value = 0.62331749999999997
$ws2.Range("E11").Value = 0.41273315999999999
$ws2.Range("E12").Value = 0.95089309
$ws2.Range("E13").Value = 0.67716052000000004
$ws2.Range("E14").Value = 0.96995677999999996
$ws2.Range("E15").Value = 0.86494046999999996
$ws2.Range("E16").Value = 1
$ws2.Range("E17").Value = 0.9908439
$ws2.Range("E18").Value = 0.69934826000000005
$ws2.Range("E19").Value = 0.55472896000000005
$ws2.Range("E20").Value = 0.21296287
$ws2.Range("E21").Value = 0.23320036999999999
$ws2.Range("E22").Value = 0.25451567000000003
$ws2.Range("E23").Value = 0.25340002
$ws2.Range("E24").Value = 0.42037773000000001
$ws2.Range("E25").Value = 0.42917559999999999
$ws2.Range("E26").Value = 0.38294892000000003
$ws2.Range("E27").Value = 0.33321145000000002
$ws2.Range("E28").Value = 0.44392799999999999
$ws2.Range("E29").Value = 0.39107617
$ws2.Range("E30").Value = 0.59904449999999998
$ws2.Range("E31").Value = 0.46894513999999998
$ws2.Range("E32").Value = 0.55383848999999996

# F3 holds its own formula; F4:F32 share the same formula pattern (relative E column)
$ws2.Range("F3").Formula = "=`$C`$3*E3^`$C`$8"
$ws2.Range("F4:F32").FormulaR1C1 = "=R3C3*RC[-1]^R8C3"

# Rows 11 & 12: clear the pre-existing empty-cell style on E/F only (G:J keep their style)
$ws2.Range("E11").Style = "Normal"
$ws2.Range("F11").Style = "Normal"
$ws2.Range("E12").Style = "Normal"
$ws2.Range("F12").Style = "Normal"

# --- Sheet "Data PGK B.Makanan (p)": column I becomes a static value (copied from the
#     β-sheet F column) instead of the old cross-sheet power formula ---
$ws3.Range("I2").Value = 422.38844218569744
$ws3.Range("I3").Value = 326.29591068314159
$ws3.Range("I4").Value = 336.61881782920665
$ws3.Range("I5").Value = 408.83820035477186
$ws3.Range("I6").Value = 417.27033754172004
$ws3.Range("I7").Value = 413.72649053101776
$ws3.Range("I8").Value = 405.69505471167224
$ws3.Range("I9").Value = 366.36966339264012
$ws3.Range("I10").Value = 301.27411744331698
$ws3.Range("I11").Value = 447.66869148048175
$ws3.Range("I12").Value = 381.06034744613902
$ws3.Range("I13").Value = 451.9052592277219
$ws3.Range("I14").Value = 427.98884959337641
$ws3.Range("I15").Value = 458.49396999999999
$ws3.Range("I16").Value = 456.4971208959895
$ws3.Range("I17").Value = 386.9348953499221
$ws3.Range("I18").Value = 346.6533327721375
$ws3.Range("I19").Value = 220.09072145762926
$ws3.Range("I20").Value = 229.77870078075213
$ws3.Range("I21").Value = 239.51596890343535
$ws3.Range("I22").Value = 239.0171974183319
$ws3.Range("I23").Value = 303.90923112770491
$ws3.Range("I24").Value = 306.91091510984853
$ws3.Range("I25").Value = 290.75450850548339
$ws3.Range("I26").Value = 272.17964144839414
$ws3.Range("I27").Value = 311.87249835541377
$ws3.Range("I28").Value = 293.66642827546781
$ws3.Range("I29").Value = 359.52902622129074
$ws3.Range("I30").Value = 320.09221825838029
$ws3.Range("I31").Value = 346.38917110053376

# --- View state: update selections to match the edited ranges ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("F3:F32").Select()

$ws3.Activate()
$ws3.Range("J2").Select()
